$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string used by A6 ("DN_SLAD" -> "DL_SLAD")
$ws.Range("A6").Value = "DL_SLAD"

# Row 3 tweaks
$ws.Range("C3").Value = 0.07095035730570878
$ws.Range("H3").Value = 0.6783968666656957

# Row 6 tweaks
$ws.Range("B6").Value = 0.9977169797572784
$ws.Range("C6").Value = 0.02668426274565401
$ws.Range("D6").Value = 0.1558647452979633
$ws.Range("E6").Value = 0.0344592836737122
$ws.Range("F6").Value = 0.9286539000248304
$ws.Range("G6").Value = 0.1567553855754553
$ws.Range("H6").Value = 0.6594255423022942
$ws.Range("I6").Value = 0.1926355703797797
